# Auto-generated edit script: refresh Seraph market-price columns (H/I/J/K/L/M/N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR per scheduled price-data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 126
$ws.Range("I11").Value = 126
$ws.Range("K11").Value = 126
$ws.Range("M11").Value = 14

$ws.Range("H74").Value = 2999.3333
$ws.Range("I74").Value = 2999.3333
$ws.Range("K74").Value = 2999.3333
$ws.Range("M74").Value = -2063.3333

$ws.Range("H77").Value = 2999.3333
$ws.Range("I77").Value = 2999.3333
$ws.Range("K77").Value = 14996.6665
$ws.Range("M77").Value = -10316.6665

$ws.Range("H132").Value = 1482.8049
$ws.Range("I132").Value = 1357.2858
$ws.Range("K132").Value = 4071.8574
$ws.Range("M132").Value = -1541.8574

$ws.Range("H135").Value = 749.3214
$ws.Range("I135").Value = 489.9524
$ws.Range("K135").Value = 4409.5716
$ws.Range("M135").Value = -1874.5716

$ws.Range("H137").Value = 2788.625
$ws.Range("I137").Value = 1559.8125
$ws.Range("K137").Value = 4679.4375
$ws.Range("M137").Value = -2129.4375

$ws.Range("H138").Value = 27459
$ws.Range("I138").Value = 16824
$ws.Range("K138").Value = 50472
$ws.Range("M138").Value = -45332

$ws.Range("H141").Value = 2647.25
$ws.Range("I141").Value = 434.8
$ws.Range("K141").Value = 1304.4
$ws.Range("M141").Value = 3875.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15558.737
$ws.Range("I32").Value = 6694.5835
$ws.Range("J32").Value = 28323.12
$ws.Range("K32").Value = 6694.5835
$ws.Range("L32").Value = 28323.12
$ws.Range("M32").Value = -6407.5835
$ws.Range("N32").Value = -28897.12

$ws.Range("H61").Value = 1405.8422
$ws.Range("I61").Value = 1405.8422
$ws.Range("K61").Value = 1405.8422
$ws.Range("M61").Value = -1193.8422

$ws.Range("H74").Value = 5102.6924
$ws.Range("I74").Value = 1982.25
$ws.Range("J74").Value = 6489.5557
$ws.Range("K74").Value = 1982.25
$ws.Range("L74").Value = 6489.5557
$ws.Range("M74").Value = -1108.25
$ws.Range("N74").Value = -8237.555700000001

$ws.Range("H77").Value = 5102.6924
$ws.Range("I77").Value = 1982.25
$ws.Range("J77").Value = 6489.5557
$ws.Range("K77").Value = 9911.25
$ws.Range("L77").Value = 32447.7785
$ws.Range("M77").Value = -5543.25
$ws.Range("N77").Value = -41183.7785

$ws.Range("H110").Value = 11719.2
$ws.Range("I110").Value = 11719.2
$ws.Range("K110").Value = 11719.2
$ws.Range("M110").Value = -9674.200000000001

$ws.Range("H132").Value = 2281.9375
$ws.Range("I132").Value = 2281.9375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6845.8125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4315.8125
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 1405.8422
$ws.Range("I136").Value = 1405.8422
$ws.Range("K136").Value = 4217.5266
$ws.Range("M136").Value = -1667.5266

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1454.5555
$ws.Range("I64").Value = 1298.75
$ws.Range("J64").Value = 1579.2
$ws.Range("K64").Value = 1298.75
$ws.Range("L64").Value = 1579.2
$ws.Range("M64").Value = -1073.75
$ws.Range("N64").Value = -2029.2

$ws.Range("H67").Value = 1454.5555
$ws.Range("I67").Value = 1298.75
$ws.Range("J67").Value = 1579.2
$ws.Range("K67").Value = 1298.75
$ws.Range("L67").Value = 1579.2
$ws.Range("M67").Value = -518.75
$ws.Range("N67").Value = -3139.2

$ws.Range("H107").Value = 2690.2727
$ws.Range("I107").Value = 2232.5557
$ws.Range("K107").Value = 2232.5557
$ws.Range("M107").Value = -312.5556999999999

$ws.Range("H134").Value = 1508.7142
$ws.Range("I134").Value = 760.3333
$ws.Range("J134").Value = 5999
$ws.Range("K134").Value = 2280.9999
$ws.Range("L134").Value = 17997
$ws.Range("M134").Value = 254.0001000000002
$ws.Range("N134").Value = -23067

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4458.7334
$ws.Range("I31").Value = 1680.5
$ws.Range("J31").Value = 5469
$ws.Range("K31").Value = 1680.5
$ws.Range("L31").Value = 5469
$ws.Range("M31").Value = -1385.5
$ws.Range("N31").Value = -6059

$ws.Range("H34").Value = 4458.7334
$ws.Range("I34").Value = 1680.5
$ws.Range("J34").Value = 5469
$ws.Range("K34").Value = 1680.5
$ws.Range("L34").Value = 5469
$ws.Range("M34").Value = -1478.5
$ws.Range("N34").Value = -5873

$ws.Range("H62").Value = 44396.3
$ws.Range("I62").Value = 4895.4
$ws.Range("J62").Value = 83897.2
$ws.Range("K62").Value = 4895.4
$ws.Range("L62").Value = 83897.2
$ws.Range("M62").Value = -4271.4
$ws.Range("N62").Value = -85145.2

$ws.Range("H65").Value = 44396.3
$ws.Range("I65").Value = 4895.4
$ws.Range("J65").Value = 83897.2
$ws.Range("K65").Value = 24477
$ws.Range("L65").Value = 419486
$ws.Range("M65").Value = -21357
$ws.Range("N65").Value = -425726

$ws.Range("H134").Value = 5112.5625
$ws.Range("I134").Value = 3778.5
$ws.Range("J134").Value = 7336
$ws.Range("K134").Value = 11335.5
$ws.Range("L134").Value = 22008
$ws.Range("M134").Value = -8800.5
$ws.Range("N134").Value = -27078

$ws.Range("H141").Value = 120000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 120000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 120000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -130360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 433.16666
$ws.Range("I5").Value = 399.75
$ws.Range("K5").Value = 1199.25
$ws.Range("M5").Value = -1087.25

$ws.Range("H99").Value = 1213.8
$ws.Range("I99").Value = 1213.8
$ws.Range("K99").Value = 3641.4
$ws.Range("M99").Value = -1395.4

$ws.Range("H106").Value = 10029
$ws.Range("J106").Value = 10029
$ws.Range("L106").Value = 30087
$ws.Range("N106").Value = -31979

$ws.Range("H113").Value = 1795.5385
$ws.Range("J113").Value = 1536.9166
$ws.Range("L113").Value = 4610.7498
$ws.Range("N113").Value = -8950.7498

$ws.Range("H131").Value = 1457.3572
$ws.Range("I131").Value = 777.3333
$ws.Range("J131").Value = 1967.375
$ws.Range("K131").Value = 2331.9999
$ws.Range("L131").Value = 5902.125
$ws.Range("M131").Value = 2708.0001
$ws.Range("N131").Value = -15982.125

$ws.Range("H134").Value = 702.3
$ws.Range("I134").Value = 702.3
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2106.9
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 2963.1
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 433.16666
$ws.Range("I135").Value = 399.75
$ws.Range("K135").Value = 3597.75
$ws.Range("M135").Value = -1062.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2210.25
$ws.Range("J97").Value = 2193.6
$ws.Range("L97").Value = 2193.6
$ws.Range("N97").Value = -3185.6

$ws.Range("H113").Value = 4999.857
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 4999.857
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 4999.857
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9339.857

$ws.Range("H122").Value = 690258.75
$ws.Range("I122").Value = 92741.55
$ws.Range("K122").Value = 278224.65
$ws.Range("M122").Value = -275774.65

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7282.7144
$ws.Range("I61").Value = 9329
$ws.Range("K61").Value = 9329
$ws.Range("M61").Value = -9127

$ws.Range("H68").Value = 3041
$ws.Range("I68").Value = 3051.25
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 3051.25
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -2302.25
$ws.Range("N68").Value = -4498

$ws.Range("H71").Value = 3041
$ws.Range("I71").Value = 3051.25
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 15256.25
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -11512.25
$ws.Range("N71").Value = -22488

$ws.Range("H104").Value = 31956.334
$ws.Range("J104").Value = 31956.334
$ws.Range("L104").Value = 31956.334
$ws.Range("N104").Value = -38944.334

$ws.Range("H113").Value = 7282.7144
$ws.Range("I113").Value = 9329
$ws.Range("K113").Value = 9329
$ws.Range("M113").Value = -7159

$ws.Range("H122").Value = 4998.8
$ws.Range("I122").Value = 4998.75
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 14996.25
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -12546.25
$ws.Range("N122").Value = -19897

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7520.294
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 7677.8125
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 7677.8125
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -8925.8125

$ws.Range("H65").Value = 7520.294
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 7677.8125
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 38389.0625
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -44629.0625

$ws.Range("H113").Value = 2231.2
$ws.Range("J113").Value = 4903
$ws.Range("L113").Value = 14709
$ws.Range("N113").Value = -19049
